$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record Burndown / Quickfire meeting values for Day 4 (column G)
$ws.Range("G12").Value = 2
$ws.Range("G14").Value = 4
$ws.Range("G21").Value = 1

# Recalculate formulas (H3 burndown total, chart caches, etc.)
$excel.Calculate()

# Move the active selection to reflect where the user left off editing
$ws.Range("G22").Select()
